# Applies the "Updated cryptos list" data refresh described by the commit diff.
# All edits are literal value replacements (prices / % changes / a few row
# reshuffles) within the existing Sheet1 table (rows 2-51, columns B-E).
# Column D is forced to Text format before writing so that numeric-looking
# price strings (e.g. "19.60", "62.80", "1.00") keep their exact text
# representation instead of being auto-converted to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.924.38"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.28"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.68"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.52"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.651.53"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.542"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0755"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.80"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.916.56"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.85"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.24"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.94"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0497"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.136.69"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.549"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.24"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.45"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.769.74"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.27"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0525"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.46"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("E51").Value = "  +1.81%  "
